$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 166667230
$ws.Range("I28").Value = 166667230
$ws.Range("K28").Value = 166667230
$ws.Range("M28").Value = -166666745
$ws.Range("H41").Value = 102
$ws.Range("I41").Value = 102
$ws.Range("K41").Value = 102
$ws.Range("M41").Value = 338
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("N47").ClearContents()
$ws.Range("H51").Value = 15294.263
$ws.Range("I51").Value = 9750
$ws.Range("K51").Value = 9750
$ws.Range("M51").Value = -9266
$ws.Range("H106").Value = 32277044
$ws.Range("I106").Value = 38477744
$ws.Range("K106").Value = 38477744
$ws.Range("M106").Value = -38477113
$ws.Range("H129").Value = 2492
$ws.Range("I129").Value = 3349.25
$ws.Range("J129").Value = 2111
$ws.Range("K129").Value = 10047.75
$ws.Range("L129").Value = 6333
$ws.Range("M129").Value = -5047.75
$ws.Range("N129").Value = -16333
$ws.Range("H137").Value = 2768.7646
$ws.Range("I137").Value = 4006.923
$ws.Range("J137").Value = 2002.2858
$ws.Range("K137").Value = 12020.769
$ws.Range("L137").Value = 6006.857400000001
$ws.Range("M137").Value = -9470.769
$ws.Range("N137").Value = -11106.8574
$ws.Range("H138").Value = 3185.1592
$ws.Range("I138").Value = 2632.7144
$ws.Range("J138").Value = 3689.5652
$ws.Range("K138").Value = 7898.1432
$ws.Range("L138").Value = 11068.6956
$ws.Range("M138").Value = -2758.1432
$ws.Range("N138").Value = -21348.6956

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 110
$ws.Range("I5").Value = 110
$ws.Range("K5").Value = 110
$ws.Range("M5").Value = 2
$ws.Range("H32").Value = 338206.4
$ws.Range("I32").Value = 4190.636
$ws.Range("K32").Value = 4190.636
$ws.Range("M32").Value = -3903.636
$ws.Range("H74").Value = 2542.3076
$ws.Range("I74").Value = 2235
$ws.Range("K74").Value = 2235
$ws.Range("M74").Value = -1361
$ws.Range("H77").Value = 2542.3076
$ws.Range("I77").Value = 2235
$ws.Range("K77").Value = 11175
$ws.Range("M77").Value = -6807

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 110
$ws.Range("I4").Value = 110
$ws.Range("K4").Value = 110
$ws.Range("M4").Value = 5
$ws.Range("H12").Value = 1871.1818
$ws.Range("I12").Value = 2370
$ws.Range("J12").Value = 1455.5
$ws.Range("K12").Value = 2370
$ws.Range("L12").Value = 1455.5
$ws.Range("M12").Value = -2202
$ws.Range("N12").Value = -1791.5
$ws.Range("H33").Value = 25000
$ws.Range("I33").Value = 25000
$ws.Range("K33").Value = 25000
$ws.Range("M33").Value = -24664
$ws.Range("H134").Value = 4764.857
$ws.Range("I134").Value = 4633.3335
$ws.Range("K134").Value = 13900.0005
$ws.Range("M134").Value = -11365.0005

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4714.5
$ws.Range("I31").Value = 4482
$ws.Range("K31").Value = 4482
$ws.Range("M31").Value = -4187
$ws.Range("H34").Value = 4714.5
$ws.Range("I34").Value = 4482
$ws.Range("K34").Value = 4482
$ws.Range("M34").Value = -4280
$ws.Range("H62").Value = 63466.145
$ws.Range("I62").Value = 7010.8
$ws.Range("K62").Value = 7010.8
$ws.Range("M62").Value = -6386.8
$ws.Range("H65").Value = 63466.145
$ws.Range("I65").Value = 7010.8
$ws.Range("K65").Value = 35054
$ws.Range("M65").Value = -31934
$ws.Range("H99").Value = 3449.8
$ws.Range("J99").Value = 2150
$ws.Range("L99").Value = 2150
$ws.Range("N99").Value = -5146
$ws.Range("H126").Value = 3449.8
$ws.Range("J126").Value = 2150
$ws.Range("L126").Value = 6450
$ws.Range("N126").Value = -11390
$ws.Range("H132").Value = 2405.875
$ws.Range("I132").Value = 2749.4
$ws.Range("K132").Value = 8248.200000000001
$ws.Range("M132").Value = -5718.200000000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 73.888885
$ws.Range("J2").Value = 89
$ws.Range("L2").Value = 534
$ws.Range("N2").Value = -760
$ws.Range("H25").Value = 234.75
$ws.Range("I25").Value = 146.33333
$ws.Range("K25").Value = 438.99999
$ws.Range("M25").Value = -269.99999
$ws.Range("H30").Value = 234.75
$ws.Range("I30").Value = 146.33333
$ws.Range("K30").Value = 438.99999
$ws.Range("M30").Value = -336.99999
$ws.Range("H40").Value = 61.88889
$ws.Range("I40").Value = 15.8
$ws.Range("J40").Value = 119.5
$ws.Range("K40").Value = 63.2
$ws.Range("L40").Value = 478
$ws.Range("M40").Value = 5.799999999999997
$ws.Range("N40").Value = -616
$ws.Range("H46").Value = 1252470
$ws.Range("J46").Value = 2002663
$ws.Range("L46").Value = 6007989
$ws.Range("N46").Value = -6008171
$ws.Range("H69").Value = 3928
$ws.Range("J69").Value = 4037.5
$ws.Range("L69").Value = 12112.5
$ws.Range("N69").Value = -13734.5
$ws.Range("H72").Value = 3928
$ws.Range("J72").Value = 4037.5
$ws.Range("L72").Value = 36337.5
$ws.Range("N72").Value = -44449.5
$ws.Range("H120").Value = 12815
$ws.Range("I120").Value = 630
$ws.Range("K120").Value = 1890
$ws.Range("M120").Value = 2948
$ws.Range("H132").Value = 3687.6155
$ws.Range("I132").Value = 2158.6667
$ws.Range("J132").Value = 4998.143
$ws.Range("K132").Value = 19428.0003
$ws.Range("L132").Value = 44983.287
$ws.Range("M132").Value = -16898.0003
$ws.Range("N132").Value = -50043.287

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 25001.5
$ws.Range("I5").Value = 4
$ws.Range("K5").Value = 4
$ws.Range("M5").Value = 108
$ws.Range("H11").Value = 5001500.5
$ws.Range("I11").Value = 7002100
$ws.Range("J11").Value = 1
$ws.Range("K11").Value = 7002100
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = -7001961
$ws.Range("N11").Value = -279
$ws.Range("H29").Value = 1000
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H33").Value = 11000
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 11000
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 11000
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -11504
$ws.Range("H36").Value = 4359.5713
$ws.Range("J36").Value = 4750
$ws.Range("L36").Value = 4750
$ws.Range("N36").Value = -5720
$ws.Range("H41").Value = 5000
$ws.Range("I41").Value = 5000
$ws.Range("K41").Value = 5000
$ws.Range("M41").Value = -4645
$ws.Range("H46").Value = 15420.417
$ws.Range("I46").Value = 5842.5
$ws.Range("J46").Value = 24998.334
$ws.Range("K46").Value = 5842.5
$ws.Range("L46").Value = 24998.334
$ws.Range("M46").Value = -5686.5
$ws.Range("N46").Value = -25310.334
$ws.Range("H132").Value = 3179.8
$ws.Range("I132").Value = 2500
$ws.Range("K132").Value = 7500
$ws.Range("M132").Value = -4970

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1999
$ws.Range("J22").Value = 1999
$ws.Range("L22").Value = 1999
$ws.Range("N22").Value = -2589
$ws.Range("H23").Value = 5000
$ws.Range("I23").Value = 5000
$ws.Range("K23").Value = 5000
$ws.Range("M23").Value = -4770
$ws.Range("H27").Value = 1999
$ws.Range("J27").Value = 1999
$ws.Range("L27").Value = 1999
$ws.Range("N27").Value = -2213
$ws.Range("H45").Value = 30000
$ws.Range("I45").Value = 30000
$ws.Range("K45").Value = 30000
$ws.Range("M45").Value = -29593

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 500
$ws.Range("J2").Value = 500
$ws.Range("L2").Value = 500
$ws.Range("N2").Value = -724
$ws.Range("H3").Value = 33747.7
$ws.Range("I3").Value = 39684.625
$ws.Range("K3").Value = 39684.625
$ws.Range("M3").Value = -39570.625
$ws.Range("H132").Value = 4863.0454
$ws.Range("I132").Value = 2525.5881
$ws.Range("K132").Value = 7576.7643
$ws.Range("M132").Value = -5046.7643
$ws.Range("H136").Value = 11356.333
$ws.Range("I136").Value = 10897
$ws.Range("K136").Value = 32691
$ws.Range("M136").Value = -30141
